# Update the test data:
#  - Remove the old "pageNavigation" sheet (sheetId 3) entirely.
#  - Rename "pageNavigation (2)" (sheetId 4) to "pageNavigation" so it
#    becomes the sole/second "pageNavigation" sheet, keeping its data
#    (rows for Appliances, TV & Audio, Computing, Gaming, Cameras,
#    Phones, Smart Tech, Home & Outdoor).
#  - Make it the active/selected sheet with C14 selected.

$wb = $excel.ActiveWorkbook

# Delete the original "pageNavigation" sheet (sheetId=3) - its rows are
# superseded by the richer "pageNavigation (2)" sheet's data.
$wb.Worksheets.Item("pageNavigation").Delete()

# Rename the remaining "pageNavigation (2)" sheet to "pageNavigation".
$ws = $wb.Worksheets.Item("pageNavigation (2)")
$ws.Name = "pageNavigation"

# Make it the active sheet with C14 selected (matches the saved view state).
$ws.Activate()
$ws.Range("C14").Select()
